$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-like pattern: force text number format on "Price" column cells
# before assigning values that look numeric, so Excel keeps them as text
# (matching the original inlineStr / shared-string cell type).

# --- Row swap: Solana (row 11) <-> WrappedEther (row 12) ---
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.975.04"
$ws.Range("E11").Value = "  +5.35%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.61"
$ws.Range("E12").Value = "  -0.94%  "

# --- Row swap: ImmutableX (row 32) <-> HuobiToken (row 33) ---
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.986"
$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7463"
$ws.Range("E33").Value = "  -1.86%  "

# --- Price/volume updates for other rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.509.47"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.23"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.39"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4286"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3658"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07275"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8623"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.384"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.529"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06936"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.62"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008885"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.37"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.770.30"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.138"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.81"
$ws.Range("E23").Value = "  +4.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.136.44"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.991"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.30"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.87"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.080"
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.17"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.822"
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08844"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.528"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05319"
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01932"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.791"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5065"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1658"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.284"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.45"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.68"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06488"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4670"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.605"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.54"
$ws.Range("E51").Value = "  -1.99%  "
